$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.111.19"
$ws.Range("E2").Value = "  -1.88%  "

$ws.Range("D3").Value = "3.477.38"
$ws.Range("E3").Value = "  -3.63%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'576.90"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'181.26"
$ws.Range("E6").Value = "  -4.63%  "

$ws.Range("D7").Value = "3.466.73"
$ws.Range("E7").Value = "  -3.88%  "

$ws.Range("E8").Value = "  -3.97%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +5.20%  "

$ws.Range("E11").Value = "  -3.86%  "

$ws.Range("D12").Value = "'53.26"
$ws.Range("E12").Value = "  -5.12%  "

$ws.Range("E13").Value = "  -4.28%  "

$ws.Range("E14").Value = "  -4.03%  "

$ws.Range("D15").Value = "4.024.30"
$ws.Range("E15").Value = "  -4.07%  "

$ws.Range("D16").Value = "'19.06"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("D17").Value = "69.081.18"
$ws.Range("E17").Value = "  -1.84%  "

$ws.Range("D18").Value = "3.481.90"
$ws.Range("E18").Value = "  -3.66%  "

$ws.Range("D19").Value = "'12.15"
$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").Value = "'531.78"
$ws.Range("E21").Value = "  +7.73%  "

$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -4.96%  "

$ws.Range("D23").Value = "'18.34"
$ws.Range("E23").Value = "  -5.13%  "

$ws.Range("D24").Value = "'4.44"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("D25").Value = "'4.82"
$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("D26").Value = "'95.41"
$ws.Range("E26").Value = "  -1.94%  "

$ws.Range("D27").Value = "'10.92"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("E28").Value = "  -3.08%  "

$ws.Range("E29").Value = "  -4.51%  "

$ws.Range("D30").Value = "'31.64"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("E31").Value = "  -5.58%  "

$ws.Range("D32").Value = "'12.30"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("D33").Value = "'63.41"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  -5.73%  "

$ws.Range("D35").Value = "'537.85"
$ws.Range("E35").Value = "  -7.77%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.399"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").Value = "'37.42"
$ws.Range("E38").Value = "  -4.32%  "

$ws.Range("E39").Value = "  +2.83%  "

$ws.Range("E40").Value = "  -9.24%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.133"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.324.13"
$ws.Range("E42").Value = "  +2.93%  "

$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "  -4.94%  "

$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("E45").Value = "  -8.84%  "

$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  -5.22%  "

$ws.Range("D47").Value = "'0.0430"
$ws.Range("E47").Value = "  -3.71%  "

$ws.Range("E48").Value = "  -4.35%  "

$ws.Range("D49").Value = "'8.91"
$ws.Range("E49").Value = "  -8.81%  "

$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").Value = "'136.73"
$ws.Range("E51").Value = "  -0.14%  "
